# Applies the Batterywise-analysis relabeling/update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / Row 7: Starting / Ending SoC (%) values swapped ---
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 19

# --- Row 8: label now carries units ---
$ws.Range("A8").Value = "Total distance covered (km)"

# --- Row 9: label renamed (value unchanged) ---
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# --- Row 10: label now carries units ---
$ws.Range("A10").Value = "Total SOC consumed(%)"

# --- Row 12-14: labels now carry units (values unchanged) ---
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# --- Row 15: label now carries units, value sign flipped (negative -> positive) ---
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.01250659291419451

# --- Row 16/17: Lowest/Highest Cell Voltage swapped (label & value) ---
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.46
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.095

# --- Row 18: label now carries units ---
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# --- Row 19/20: labels now carry units ---
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

# --- Row 21: label now carries units, value filled in (was blank) ---
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 9

# --- Row 22-27: labels now carry units / BMS suffix (values unchanged) ---
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# --- Row 28/29: lowest/highest cell temp swapped (label only, values unchanged) ---
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

# --- Row 30: label now carries units ---
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31: was "Maximum BMS Temperature in C"/60 -> now "Battery Voltage(V)"/54 ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 54

# --- Row 32: was "Battery Voltage"/5.4 -> now "Total energy charged(kWh)"/1.738524735 ---
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.738524735

# --- Row 33: was "Total energy charged in kWh"/0.1738524735 -> now "Electricity consumption units(kW)"/1.276224993393235e-07 ---
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 1.276224993393235 / 10000000

# --- Row 34: was "Electricity consumption units in kW"/1.276224993393235e-08 -> now "Idling time percentage"/16.95508798552543 ---
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 16.95508798552543

# --- Row 35: was "Idling time percentage"/16.20801354072431 -> now "Time spent in 0-10 km/h"/3.18673942860478 ---
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 3.18673942860478

# --- Row 36: was "Time spent in 0-10 km/h"/19.39475296932909 -> now "Time spent in 10-20 km/h"/3.318061108354977 ---
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.318061108354977

# --- Row 37: was "Time spent in 10-20 km/h"/3.323897627454986 -> now "Time spent in 20-30 km/h"/2.967869962354452 ---
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 2.967869962354452

# --- Row 38: was "Time spent in 20-30 km/h"/2.973706481454461 -> now "Time spent in 30-40 km/h"/14.21776052762133 ---
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 14.21776052762133

# --- Row 39: was "Time spent in 30-40 km/h"/14.21776052762133 -> now "Time spent in 40-50 km/h"/9.656520850964485 ---
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 9.656520850964485

# --- Row 40: was "Time spent in 40-50 km/h"/9.67403040826451 -> now "Time spent in 50-60 km/h"/42.17468701666326 ---
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 42.17468701666326

# --- Row 41: was "Time spent in 50-60 km/h"/42.19511483351329 -> now "Time spent in 60-70 km/h"/7.409460997461115 ---
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 7.409460997461115

# --- Row 42: was "Time spent in 60-70 km/h"/7.47366270756121 -> now "Time spent in 70-80 km/h"/0 ---
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0

# --- Row 43: new row ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
